# Updates the "想去人数" (F column) counts across the four worksheets
# to match the values published in the new gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$ExhibitionValues = @{
    7  = 9620
    9  = 723
    10 = 2200
    12 = 1667
    13 = 2793
    14 = 148
    15 = 4175
    16 = 354
    17 = 177
    18 = 137
    19 = 531
    23 = 90
    25 = 4085
    26 = 6
    27 = 3548
    28 = 1112
    30 = 512
    32 = 77
    33 = 366
    34 = 449
    35 = 349
}
foreach ($row in $ExhibitionValues.Keys) {
    $wsExhibition.Range("F$row").Value = $ExhibitionValues[$row]
}

# 演出 (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$PerformanceValues = @{
    3 = 25
    5 = 27
}
foreach ($row in $PerformanceValues.Keys) {
    $wsPerformance.Range("F$row").Value = $PerformanceValues[$row]
}

# 本地生活 (Local Life)
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$LocalLifeValues = @{
    2 = 204
    3 = 1014
}
foreach ($row in $LocalLifeValues.Keys) {
    $wsLocalLife.Range("F$row").Value = $LocalLifeValues[$row]
}

# 全部类型 (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$AllValues = @{
    2  = 204
    4  = 1014
    9  = 9620
    11 = 723
    12 = 2200
    14 = 1667
    16 = 2793
    17 = 148
    18 = 4175
    19 = 354
    20 = 177
    21 = 137
    22 = 531
    25 = 25
    27 = 90
    29 = 4085
    30 = 6
    31 = 3548
    32 = 1112
    34 = 512
    36 = 77
    37 = 366
    38 = 449
    39 = 349
    41 = 27
}
foreach ($row in $AllValues.Keys) {
    $wsAll.Range("F$row").Value = $AllValues[$row]
}
